# Refresh the crypto price/volume table with the latest scraped values.
# Prices in column D that Excel would otherwise auto-parse as numbers are
# written with a leading apostrophe so they stay as text (matching the
# inlineStr cells in the source data, e.g. "1.00" must not become 1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.670.06"
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("D3").Value = "2.606.87"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'594.73"
$ws.Range("E5").Value = "  -2.03%  "
$ws.Range("D6").Value = "'151.08"
$ws.Range("E6").Value = "  +3.68%  "
$ws.Range("D8").Value = "'0.587"
$ws.Range("E8").Value = "  +0.60%  "
$ws.Range("D9").Value = "'0.109"
$ws.Range("E9").Value = "  +1.53%  "
$ws.Range("D10").Value = "'5.67"
$ws.Range("E10").Value = "  +3.02%  "
$ws.Range("D11").Value = "'0.385"
$ws.Range("E11").Value = "  +3.49%  "
$ws.Range("E12").Value = "  -0.83%  "
$ws.Range("D13").Value = "'27.82"
$ws.Range("E13").Value = "  +2.43%  "
$ws.Range("D14").Value = "3.079.11"
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("D15").Value = "63.530.36"
$ws.Range("E15").Value = "  +0.83%  "
$ws.Range("D16").Value = "'0.0000154"
$ws.Range("E16").Value = "  +5.52%  "
$ws.Range("D17").Value = "2.629.30"
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("D18").Value = "'12.40"
$ws.Range("E18").Value = "  +7.97%  "
$ws.Range("E19").Value = "  +4.89%  "
$ws.Range("D20").Value = "'348.37"
$ws.Range("E20").Value = "  +2.00%  "
$ws.Range("E21").Value = "  +0.39%  "
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").Value = "'67.62"
$ws.Range("E23").Value = "  +2.48%  "
$ws.Range("E24").Value = "  +7.32%  "
$ws.Range("D25").Value = "'9.35"
$ws.Range("E25").Value = "  +3.78%  "
$ws.Range("D26").Value = "'1.69"
$ws.Range("E26").Value = "  +0.33%  "
$ws.Range("D27").Value = "'555.05"
$ws.Range("E27").Value = "  +2.22%  "
$ws.Range("D28").Value = "'8.02"
$ws.Range("E28").Value = "  +2.58%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("E31").Value = "  +1.60%  "
$ws.Range("D32").Value = "0.0₃0850"
$ws.Range("E32").Value = "  +1.20%  "
$ws.Range("D33").Value = "'1.75"
$ws.Range("E33").Value = "  +0.69%  "
$ws.Range("D34").Value = "'5.23"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").Value = "'166.88"
$ws.Range("E35").Value = "  -1.01%  "
$ws.Range("D36").Value = "'0.415"
$ws.Range("E36").Value = "  +3.34%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").Value = "'19.59"
$ws.Range("E38").Value = "  +3.76%  "
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").Value = "'166.74"
$ws.Range("E41").Value = "  +0.84%  "
$ws.Range("D42").Value = "'39.72"
$ws.Range("E42").Value = "  +0.23%  "
$ws.Range("D43").Value = "'3.96"
$ws.Range("E43").Value = "  +5.75%  "
$ws.Range("D44").Value = "'0.0588"
$ws.Range("E44").Value = "  +4.59%  "
$ws.Range("D45").Value = "'21.96"
$ws.Range("E45").Value = "  +1.09%  "
$ws.Range("D46").Value = "'0.633"
$ws.Range("E46").Value = "  +1.76%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0252"
$ws.Range("E47").Value = "  +4.39%  "
$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").Value = "'2.03"
$ws.Range("E48").Value = "  +5.26%  "
$ws.Range("D49").Value = "'0.0966"
$ws.Range("E49").Value = "  +1.26%  "
$ws.Range("D50").Value = "'19.22"
$ws.Range("E50").Value = "  +3.69%  "
$ws.Range("E51").Value = "  +22.17%  "
